# Adds the new "SE2_CI" / "SE2_CNI" columns (retention figures) to the
# historical table on "Hoja1", as new columns AH and AI, with one pair of
# figures per existing territory row. This mirrors the underlying report
# export, which rewrote the AD:AG ("SE53_CI".."SE1_CNI") figures as text
# at the same time the two new columns were appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing AD:AG ("SE53_CI", "SE53_CNI", "SE1_CI", "SE1_CNI") data, rewritten
# as text by the export that introduced the new columns (values unchanged).
$existingData = @(
    @{ Row = 2;  AD = "2209.0";  AE = "1474.0"; AF = "2304.0";  AG = "1989.0" },
    @{ Row = 3;  AD = "2306.0";  AE = "991.0";  AF = "2420.0";  AG = "1382.0" },
    @{ Row = 4;  AD = "2284.0";  AE = "1170.0"; AF = "2723.0";  AG = "1283.0" },
    @{ Row = 5;  AD = "1073.0";  AE = "482.0";  AF = "994.0";   AG = "484.0" },
    @{ Row = 6;  AD = "1447.0";  AE = "708.0";  AF = "1572.0";  AG = "593.0" },
    @{ Row = 7;  AD = "682.0";   AE = "353.0";  AF = "706.0";   AG = "356.0" },
    @{ Row = 8;  AD = "862.0";   AE = "782.0";  AF = "1035.0";  AG = "762.0" },
    @{ Row = 9;  AD = "245.0";   AE = "157.0";  AF = "306.0";   AG = "222.0" },
    @{ Row = 10; AD = "485.0";   AE = "257.0";  AF = "822.0";   AG = "526.0" },
    @{ Row = 11; AD = "11593.0"; AE = "6374.0"; AF = "12882.0"; AG = "7597.0" }
)

$ws.Range("AD2:AG11").NumberFormat = "@"
foreach ($entry in $existingData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 30).Value = $entry.AD  # column AD
    $ws.Cells.Item($r, 31).Value = $entry.AE  # column AE
    $ws.Cells.Item($r, 32).Value = $entry.AF  # column AF
    $ws.Cells.Item($r, 33).Value = $entry.AG  # column AG
}

# Header row (row 1) for the two new columns.
$ws.Range("AH1").Value = "SE2_CI"
$ws.Range("AI1").Value = "SE2_CNI"

# Data rows 2-11: new SE2_CI / SE2_CNI figures per territory, in the same
# row order as the existing table (Barcelona Ciutat ... Totals).
$newData = @(
    @{ Row = 2;  CI = 2779.0;  CNI = 3381.0 },
    @{ Row = 3;  CI = 2921.0;  CNI = 2865.0 },
    @{ Row = 4;  CI = 3578.0;  CNI = 2746.0 },
    @{ Row = 5;  CI = 1326.0;  CNI = 779.0 },
    @{ Row = 6;  CI = 1770.0;  CNI = 736.0 },
    @{ Row = 7;  CI = 863.0;   CNI = 562.0 },
    @{ Row = 8;  CI = 1375.0;  CNI = 1021.0 },
    @{ Row = 9;  CI = 511.0;   CNI = 435.0 },
    @{ Row = 10; CI = 797.0;   CNI = 729.0 },
    @{ Row = 11; CI = 15920.0; CNI = 13254.0 }
)

foreach ($entry in $newData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 34).Value = $entry.CI   # column AH
    $ws.Cells.Item($r, 35).Value = $entry.CNI  # column AI
}
